$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the "_old" / "_new" header labels to "_FV2210" / "_FV2304".
#    (Row 1 holds the column headers; "diff" in K1 is left untouched.)
# ---------------------------------------------------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2210"
    "B1" = "Segmentgruppe_FV2210"
    "C1" = "Segment_FV2210"
    "D1" = "Datenelement_FV2210"
    "E1" = "Segment ID_FV2210"
    "F1" = "Code_FV2210"
    "G1" = "Qualifier_FV2210"
    "H1" = "Beschreibung_FV2210"
    "I1" = "Bedingungsausdruck_FV2210"
    "J1" = "Bedingung_FV2210"
    "L1" = "Segmentname_FV2304"
    "M1" = "Segmentgruppe_FV2304"
    "N1" = "Segment_FV2304"
    "O1" = "Datenelement_FV2304"
    "P1" = "Segment ID_FV2304"
    "Q1" = "Code_FV2304"
    "R1" = "Qualifier_FV2304"
    "S1" = "Beschreibung_FV2304"
    "T1" = "Bedingungsausdruck_FV2304"
    "U1" = "Bedingung_FV2304"
}
foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# ---------------------------------------------------------------------------
# 2. Turn A1:U67 into a real Excel Table (ListObject), named Table1, with
#    an AutoFilter - matching the new xl/tables/table1.xml part.
#    The header row already carries explicit formatting (bold/fill/border);
#    stash it in a scratch row and restore it verbatim after the table is
#    created so no extra dxf / cell style gets synthesized for the table
#    header.
# ---------------------------------------------------------------------------
$headerRange  = $ws.Range("A1:U1")
$scratchRow   = 70
$scratchRange = $ws.Range("A" + $scratchRow + ":U" + $scratchRow)

$headerRange.Copy()
$scratchRange.PasteSpecial(-4122)   # xlPasteFormats

$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)    # xlPasteFormats
$ws.Rows($scratchRow).Delete()

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
